# Update the "想去人数" (interest count) figures that changed between the
# two generated-output snapshots.  The same underlying event list is
# duplicated across the "展览" sheet (Worksheets item 1) and the
# "全部类型" sheet (Worksheets item 4); both need the refreshed counts.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)   # 展览
$sheet4 = $wb.Worksheets.Item(4)   # 全部类型

# Row -> new value for the "展览" sheet
$updates1 = @{
    5  = 15556
    9  = 15396
    11 = 8994
    13 = 7
    21 = 549
    25 = 1110
    28 = 81
    32 = 60
    33 = 39
    34 = 248
    38 = 5517
}

foreach ($row in $updates1.Keys) {
    $sheet1.Range("F$row").Value = $updates1[$row]
}

# Row -> new value for the "全部类型" sheet
$updates4 = @{
    5  = 15556
    9  = 15396
    11 = 8994
    13 = 7
    21 = 549
    25 = 1110
    28 = 81
    34 = 60
    35 = 39
    36 = 248
    40 = 5517
}

foreach ($row in $updates4.Keys) {
    $sheet4.Range("F$row").Value = $updates4[$row]
}
